# The deck currently carries the "Integral" / "Red Violet" theme on
# ppt/theme/theme1.xml (the theme used by the slide master, i.e. the
# slides) and the stock "Office Theme" on ppt/theme/theme2.xml (the
# theme used by the notes master). The edit re-colours the slide
# theme so it matches the stock "Office Theme" palette (swapping the
# two themes' colour schemes).
#
# PowerPoint exposes the 12 theme colour slots (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) through Slide.ThemeColorScheme, in that
# fixed order, as ColorFormat objects with a settable .RGB (a VBA-style
# 0x00BBGGRR packed Long).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeColor {
    param(
        [int]$Index,
        [string]$Hex
    )
    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)
    $packed = $r + ($g * 256) + ($b * 65536)
    $tcs.Item($Index).RGB = $packed
}

# Target palette: the stock Office theme colours.
Set-ThemeColor 1  "000000"   # dk1
Set-ThemeColor 2  "FFFFFF"   # lt1
Set-ThemeColor 3  "44546A"   # dk2
Set-ThemeColor 4  "E7E6E6"   # lt2
Set-ThemeColor 5  "5B9BD5"   # accent1
Set-ThemeColor 6  "ED7D31"   # accent2
Set-ThemeColor 7  "A5A5A5"   # accent3
Set-ThemeColor 8  "FFC000"   # accent4
Set-ThemeColor 9  "4472C4"   # accent5
Set-ThemeColor 10 "70AD47"   # accent6
Set-ThemeColor 11 "0563C1"   # hlink
Set-ThemeColor 12 "954F72"   # folHlink
